$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DMD")

# ---------------------------------------------------------------------
# New "TimeSlice" block, rows 24-34.
# Write order below is deliberate: it controls the order in which new
# shared-string entries are appended (the engine appends new strings in
# first-use order), so that it matches the target file's table.
# ---------------------------------------------------------------------

# Row 24: "~FI_T:" section header - same look as B15 ("~FI_T:DEMAND")
$ws.Range("B15").Copy() | Out-Null
$ws.Range("B24").PasteSpecial(-4122) | Out-Null
$ws.Range("B24").Value = "~FI_T:"

# Row 25: column headers - same look as row 16 (CommName / year header row)
$ws.Range("B16").Copy() | Out-Null
$ws.Range("B25").PasteSpecial(-4122) | Out-Null
$ws.Range("B25").Value = "CommName"

$ws.Range("C16").Copy() | Out-Null
$ws.Range("C25").PasteSpecial(-4122) | Out-Null

$ws.Range("D25").Value = "COM_FR"
$ws.Range("D25").Interior.Color = 12713983
$ws.Range("D25").Font.Name = "Arial"
$ws.Range("D25").Font.Size = 10
$ws.Range("D25").HorizontalAlignment = -4131
$ws.Range("D25").VerticalAlignment = -4108

# Row 26: "\I: Demand Commodity Name" label, like row 17
$ws.Range("B17").Copy() | Out-Null
$ws.Range("B26").PasteSpecial(-4122) | Out-Null
$ws.Range("B26").Value = "\I: Demand Commodity Name"

$ws.Range("C26:D26").Interior.Color = 10079487
$ws.Range("C26:D26").VerticalAlignment = -4108
$ws.Range("C26:D26").WrapText = $true
$ws.Range("C26:D26").Borders.Item(9).LineStyle = 1
$ws.Range("C26:D26").Borders.Item(9).Weight = -4138
$ws.Rows.Item(26).RowHeight = 26.25

# Rows 27-34: one row per TimeSlice, alternating banded fill, matching
# the existing "ELEC_HV / value" row (row 18) for columns B and D.
$tsNames = @("1S1W1D","1S1W2D","1S2W1D","1S2W2D","2S1W1D","2S1W2D","2S2W1D","2S2W2D")
$tsValues = @(0.15,0.1,0.15,0.1,0.15,0.1,0.15,0.1)

for ($i = 0; $i -lt 8; $i++) {
    $r = 27 + $i
    $rowStr = [string]$r

    $ws.Range("B18").Copy() | Out-Null
    $ws.Range("B" + $rowStr).PasteSpecial(-4122) | Out-Null
    $ws.Range("B" + $rowStr).Formula = "=SEC_Comm!C8"

    $cCell = $ws.Range("C" + $rowStr)
    $cCell.Value = $tsNames[$i]
    $cCell.Font.Name = "Arial"
    $cCell.Font.Size = 10
    if (($i % 2) -eq 0) {
        $cCell.Interior.Color = 16777215
    } else {
        $cCell.Interior.Color = 15921906
    }

    $ws.Range("E18").Copy() | Out-Null
    $ws.Range("D" + $rowStr).PasteSpecial(-4122) | Out-Null
    $ws.Range("D" + $rowStr).Value = $tsValues[$i]
}

# Row 34 (last TimeSlice row) gets a thick bottom border under column C.
$ws.Range("C34").Borders.Item(9).LineStyle = 1
$ws.Range("C34").Borders.Item(9).Weight = -4138
$ws.Rows.Item(34).RowHeight = 13.5

# Finally the "TimeSlice" column header text, last of the new strings.
$ws.Range("C25").Value = "TimeSlice"

# ---------------------------------------------------------------------
# Active-sheet / selection bookkeeping: DMD becomes the selected tab.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("J25").Select()
